$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns are treated as text so that
# numeric-looking strings (e.g. "0.999") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2: update D2, E2
$ws.Range("D2").Value = "66.765.96"
$ws.Range("E2").Value = "  +0.45%  "

# Row 3: update D3, E3
$ws.Range("D3").Value = "3.829.24"
$ws.Range("E3").Value = "  +3.82%  "

# Row 4: update E4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5: update D5, E5
$ws.Range("D5").Value = "411.27"
$ws.Range("E5").Value = "  -1.93%  "

# Row 6: update D6, E6
$ws.Range("D6").Value = "131.40"
$ws.Range("E6").Value = "  +1.13%  "

# Row 7: update D7, E7
$ws.Range("D7").Value = "3.821.29"
$ws.Range("E7").Value = "  +3.84%  "

# Row 8: update E8
$ws.Range("E8").Value = "  -4.32%  "

# Row 9: update D9, E9
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.08%  "

# Row 10: update E10
$ws.Range("E10").Value = "  -4.24%  "

# Row 11: update E11
$ws.Range("E11").Value = "  -6.11%  "

# Row 12: update E12
$ws.Range("E12").Value = "  -6.09%  "

# Row 13: update D13, E13
$ws.Range("D13").Value = "41.05"
$ws.Range("E13").Value = "  -4.76%  "

# Row 14: update D14, E14
$ws.Range("D14").Value = "4.445.95"
$ws.Range("E14").Value = "  +3.98%  "

# Row 15: update D15, E15
$ws.Range("D15").Value = "9.99"
$ws.Range("E15").Value = "  -5.93%  "

# Row 16: update D16, E16
$ws.Range("D16").Value = "15.32"
$ws.Range("E16").Value = "  +14.77%  "

# Row 17: update B17, C17, D17, E17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.138"
$ws.Range("E17").Value = "  -1.08%  "

# Row 18: update B18, C18, D18, E18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.824.20"
$ws.Range("E18").Value = "  +3.55%  "

# Row 19: update E19
$ws.Range("E19").Value = "  -4.68%  "

# Row 20: update D20, E20
$ws.Range("D20").Value = "67.210.57"
$ws.Range("E20").Value = "  +0.98%  "

# Row 21: update E21
$ws.Range("E21").Value = "  -5.06%  "

# Row 22: update D22, E22
$ws.Range("D22").Value = "412.99"
$ws.Range("E22").Value = "  -7.00%  "

# Row 23: update D23, E23
$ws.Range("D23").Value = "14.48"
$ws.Range("E23").Value = "  -12.08%  "

# Row 24: update D24, E24
$ws.Range("D24").Value = "85.42"
$ws.Range("E24").Value = "  -5.04%  "

# Row 25: update E25
$ws.Range("E25").Value = "  -2.86%  "

# Row 26: update B26, C26, D26, E26
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "5.79"
$ws.Range("E26").Value = "  +12.91%  "

# Row 27: update B27, C27, D27, E27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "36.57"
$ws.Range("E27").Value = "  -2.06%  "

# Row 28: update E28
$ws.Range("E28").Value = "  -6.24%  "

# Row 29: update D29, E29
$ws.Range("D29").Value = "9.44"
$ws.Range("E29").Value = "  -7.52%  "

# Row 30: update D30, E30
$ws.Range("D30").Value = "680.36"
$ws.Range("E30").Value = "  +4.27%  "

# Row 31: update B31, C31, D31, E31
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.122"
$ws.Range("E31").Value = "  -2.24%  "

# Row 32: update B32, C32, D32, E32
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "12.46"
$ws.Range("E32").Value = "  -1.98%  "

# Row 33: update E33
$ws.Range("E33").Value = "  +0.47%  "

# Row 34: update E34
$ws.Range("E34").Value = "  -1.63%  "

# Row 35: update E35
$ws.Range("E35").Value = "  -8.79%  "

# Row 36: update D36, E36
$ws.Range("D36").Value = "38.74"
$ws.Range("E36").Value = "  -7.07%  "

# Row 37: update E37
$ws.Range("E37").Value = "  -0.02%  "

# Row 38: update D38, E38
$ws.Range("D38").Value = "0.0₃0791"
$ws.Range("E38").Value = "  +8.55%  "

# Row 39: update D39, E39
$ws.Range("D39").Value = "55.12"
$ws.Range("E39").Value = "  -3.76%  "

# Row 40: update D40, E40
$ws.Range("D40").Value = "3.10"
$ws.Range("E40").Value = "  +0.08%  "

# Row 41: update D41, E41
$ws.Range("D41").Value = "0.0459"
$ws.Range("E41").Value = "  -6.96%  "

# Row 42: update D42, E42
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  +0.00%  "

# Row 43: update E43
$ws.Range("E43").Value = "  -8.67%  "

# Row 44: update D44, E44
$ws.Range("D44").Value = "149.19"
$ws.Range("E44").Value = "  +0.15%  "

# Row 45: update D45, E45
$ws.Range("D45").Value = "4.50"
$ws.Range("E45").Value = "  +3.34%  "

# Row 46: update D46, E46
$ws.Range("D46").Value = "3.32"
$ws.Range("E46").Value = "  -2.99%  "

# Row 47: update D47, E47
$ws.Range("D47").Value = "3.16"
$ws.Range("E47").Value = "  +17.38%  "

# Row 48: update D48, E48
$ws.Range("D48").Value = "26.79"
$ws.Range("E48").Value = "  -8.81%  "

# Row 49: update D49, E49
$ws.Range("D49").Value = "2.09"
$ws.Range("E49").Value = "  -1.08%  "

# Row 50: update E50
$ws.Range("E50").Value = "  -3.20%  "

# Row 51: update D51, E51
$ws.Range("D51").Value = "2.55"
$ws.Range("E51").Value = "  -4.22%  "

# Reset style to Normal so the temporary text NumberFormat does not linger
# as a persisted custom style on the cells (restores default appearance).
$ws.Range("D2:E51").Style = "Normal"
